$wb = $excel.ActiveWorkbook

# Rename the second sheet ("Tabelle2") to "hidden_SheetName"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "hidden_SheetName"

# Make the renamed sheet the active tab (moves tabSelected from sheet1 to
# sheet2 and sets workbookView activeTab to its 0-based index).
$ws2.Activate()
